# Regenerate merged AHB files
# 1. Rename the "_old"/"_new" header suffixes to the specific version tags
#    ("_FV2310" / "_FV2404") used by this merged AHB diff export.
# 2. Turn the header/data range into a native Excel table ("Table1").
# 3. Freeze the header row (row 1) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels -------------------------------------------------
$oldHeaders = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

$newHeaders = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    # columns A..J (1..10) hold the "_old" -> "_FV2310" headers
    $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}

for ($i = 0; $i -lt $newHeaders.Count; $i++) {
    # columns L..U (12..21) hold the "_new" -> "_FV2404" headers (column K is "diff")
    $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i]
}

# --- 2. Convert the used range into a table ---------------------------------------
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U54"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
